$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the two new "attachment" rows (week-2 homework additions) ---
# Labels first (column A), top to bottom.
$ws.Range("A11").Value = "[参考]厨房.zip"
$ws.Range("A12").Value = "[参考]客厅.rar"

# Then the link targets/text (column B) - B12 before B11.
$ws.Range("B12").Value = "https://github.com/MuzhiYing/vr_indoorFurnitureRoaming/blob/master/01%E9%9C%80%E6%B1%82%E8%B0%83%E7%A0%94/%E9%99%84%E5%BD%95/%5B%E5%8F%82%E8%80%83%5D%E5%AE%A2%E5%8E%85.rar"
$ws.Range("B11").Value = "https://github.com/MuzhiYing/vr_indoorFurnitureRoaming/blob/master/01%E9%9C%80%E6%B1%82%E8%B0%83%E7%A0%94/%E9%99%84%E5%BD%95/%5B%E5%8F%82%E8%80%83%5D%E5%8E%A8%E6%88%BF.zip"

# Wire up the hyperlinks themselves (B12 first, then B11 - matches authoring order).
$ws.Hyperlinks.Add($ws.Range("B12"), $ws.Range("B12").Value)
$ws.Hyperlinks.Add($ws.Range("B11"), $ws.Range("B11").Value)

# --- View state: scroll back to the top and land the selection on B11 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("B11").Select()
